$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 174096
$ws.Range("C4").Value = 164090
$ws.Range("C5").Value = 10007
$ws.Range("C8").Value = 64.42
